# Update automatico via Actualizar 02-19-2021 12-14-56
# Refresh the "Ultimo" (last updated) timestamp column (D) for each of the
# three stacked update-batches in the sheet, shifting the recorded serial
# date/time values to the newly observed ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value  = 44246.51022088
$ws.Range("D16:D29").Value = 44246.48894626158
$ws.Range("D30:D43").Value = 44246.46755512732
